# Auto-applies cell value updates to match the target diff (cryptos price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '60.349.62'
$ws.Range("D3").Value2 = '2.622.44'
$ws.Range("E3").Value2 = '  +0.67%  '
$ws.Range("E4").Value2 = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '519.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '150.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  -1.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  +0.07%  '
$ws.Range("E8").Value2 = '  -4.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '6.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  -4.75%  '
$ws.Range("E10").Value2 = '  +1.90%  '
$ws.Range("E11").Value2 = '  -0.53%  '
$ws.Range("E12").Value2 = '  -0.94%  '
$ws.Range("D13").Value2 = '3.081.78'
$ws.Range("E13").Value2 = '  +0.67%  '
$ws.Range("D14").Value2 = '60.366.50'
$ws.Range("E14").Value2 = '  -0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '21.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = '  -0.64%  '
$ws.Range("E16").Value2 = '  -0.61%  '
$ws.Range("D17").Value2 = '2.628.07'
$ws.Range("E17").Value2 = '  +0.78%  '
$ws.Range("E18").Value2 = '  -1.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '346.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  -3.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '10.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '6.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '0.994'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  -0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '60.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '0.422'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  -0.86%  '
$ws.Range("E25").Value2 = '  -0.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  +0.86%  '
$ws.Range("D27").Value2 = '0.0₃0835'
$ws.Range("E27").Value2 = '  -0.59%  '
$ws.Range("E28").Value2 = '  -2.84%  '
$ws.Range("E29").Value2 = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '6.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +2.94%  '
$ws.Range("E31").Value2 = '  +0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '19.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  -1.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '149.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  -0.84%  '
$ws.Range("E34").Value2 = '  -0.57%  '
$ws.Range("E35").Value2 = '  -2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.895'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  +0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.880'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = '  +4.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '36.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  +1.32%  '
$ws.Range("E39").Value2 = '  -1.79%  '
$ws.Range("E40").Value2 = '  -1.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '289.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  -0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.630'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  +1.53%  '
$ws.Range("E43").Value2 = '  -1.30%  '
$ws.Range("E44").Value2 = '  +0.22%  '
$ws.Range("E45").Value2 = '  -0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '19.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  -0.01%  '
$ws.Range("E48").Value2 = '  -4.06%  '
$ws.Range("E49").Value2 = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '18.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  -1.21%  '
$ws.Range("D51").Value2 = '1.965.46'
$ws.Range("E51").Value2 = '  -1.22%  '
